$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: insert 3 new rows before the old row 17 (Engine Speed row) ---
# After this, old row 17 becomes row 20.
$ws1.Rows("17:19").Insert()

# Copy formatting (cell styles/borders) from row 16 into the new rows 17-19
$ws1.Range("A16:C16").Copy() | Out-Null
$ws1.Range("A17:C17").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws1.Range("A18:C18").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws1.Range("A19:C19").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Update Engine Power value
$ws1.Range("B15").Value = 20000

# Row 16: SFOC at 100%
$ws1.Range("A16").Value = "sfoc at 100%"
$ws1.Range("C16").Value = "g/kwh"

# Row 17: SFOC at 85%
$ws1.Range("A17").Value = "sfoc at 85%"
$ws1.Range("B17").Value = 173
$ws1.Range("C17").Value = "g/kwh"

# Row 18: SFOC at 75%
$ws1.Range("A18").Value = "sfoc at 75%"
$ws1.Range("B18").Value = 176.4
$ws1.Range("C18").Value = "g/kwh"

# Row 19: SFOC at 50%
$ws1.Range("A19").Value = "sfoc at 50%"
$ws1.Range("B19").Value = 179.3
$ws1.Range("C19").Value = "g/kwh"

# Update dimension / view
$ws1.Application.ActiveWindow.ScrollRow = 3
$ws1.Range("B16").Select()
